$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.480.05'
$ws.Range('D2').ClearFormats()

$ws.Range('E2').Value = '  -0.44%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.633.36'
$ws.Range('D3').ClearFormats()

$ws.Range('E3').Value = '  -0.52%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9997'
$ws.Range('D4').ClearFormats()

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9998'
$ws.Range('D5').ClearFormats()

$ws.Range('E5').Value = '  +0.07%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '305.59'
$ws.Range('D6').ClearFormats()

$ws.Range('E6').Value = '  -0.99%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3740'
$ws.Range('D7').ClearFormats()

$ws.Range('E7').Value = '  -0.56%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3668'
$ws.Range('D8').ClearFormats()

$ws.Range('E8').Value = '  -0.23%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '51.90'
$ws.Range('D9').ClearFormats()

$ws.Range('E9').Value = '  -1.79%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08203'
$ws.Range('D10').ClearFormats()

$ws.Range('E10').Value = '  +0.04%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.232'
$ws.Range('D11').ClearFormats()

$ws.Range('E11').Value = '  -3.80%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.0000'
$ws.Range('D12').ClearFormats()

$ws.Range('E12').Value = '  +0.10%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.61'
$ws.Range('D13').ClearFormats()

$ws.Range('E13').Value = '  -2.04%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.585'
$ws.Range('D14').ClearFormats()

$ws.Range('E14').Value = '  -1.51%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001251'
$ws.Range('D15').ClearFormats()

$ws.Range('E15').Value = '  -2.77%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.288'
$ws.Range('D16').ClearFormats()

$ws.Range('E16').Value = '  -2.22%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.635.35'
$ws.Range('D17').ClearFormats()

$ws.Range('E17').Value = '  -0.24%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '94.48'
$ws.Range('D18').ClearFormats()

$ws.Range('E18').Value = '  -0.66%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06976'
$ws.Range('D19').ClearFormats()

$ws.Range('E19').Value = '  +0.68%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.81'
$ws.Range('D20').ClearFormats()

$ws.Range('E20').Value = '  -3.08%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.477'
$ws.Range('D21').ClearFormats()

$ws.Range('E21').Value = '  -1.69%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('D22').ClearFormats()

$ws.Range('E22').Value = '  +0.28%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.79'
$ws.Range('D23').ClearFormats()

$ws.Range('E23').Value = '  -1.06%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.469.10'
$ws.Range('D24').ClearFormats()

$ws.Range('E24').Value = '  -0.46%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.197'
$ws.Range('D25').ClearFormats()

$ws.Range('E25').Value = '  +3.06%  '

$ws.Range('E26').Value = '  +1.88%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.43'
$ws.Range('D27').ClearFormats()

$ws.Range('E27').Value = '  +0.28%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '149.97'
$ws.Range('D28').ClearFormats()

$ws.Range('E28').Value = '  -1.07%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.328'
$ws.Range('D29').ClearFormats()

$ws.Range('E29').Value = '  -0.31%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.90'
$ws.Range('D30').ClearFormats()

$ws.Range('E30').Value = '  -1.02%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.814.11'
$ws.Range('D31').ClearFormats()

$ws.Range('E31').Value = '  -0.51%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.276'
$ws.Range('D32').ClearFormats()

$ws.Range('E32').Value = '  -4.68%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.843'
$ws.Range('D33').ClearFormats()

$ws.Range('E33').Value = '  -0.34%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.027'
$ws.Range('D34').ClearFormats()

$ws.Range('E34').Value = '  +3.82%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '10.94'
$ws.Range('D35').ClearFormats()

$ws.Range('E35').Value = '  +4.67%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02794'
$ws.Range('D36').ClearFormats()

$ws.Range('E36').Value = '  -2.32%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2539'
$ws.Range('D37').ClearFormats()

$ws.Range('E37').Value = '  -0.73%  '

$ws.Range('B38').Value = 'InternetComputer(DFINITY)'

$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.088'
$ws.Range('D38').ClearFormats()

$ws.Range('E38').Value = '  -2.10%  '

$ws.Range('B39').Value = 'Stellar'

$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.08759'
$ws.Range('D39').ClearFormats()

$ws.Range('E39').Value = '  -1.68%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.07144'
$ws.Range('D40').ClearFormats()

$ws.Range('E40').Value = '  -3.61%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.7086'
$ws.Range('D41').ClearFormats()

$ws.Range('E41').Value = '  -0.90%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.353'
$ws.Range('D42').ClearFormats()

$ws.Range('E42').Value = '  -2.59%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.27'
$ws.Range('D43').ClearFormats()

$ws.Range('E43').Value = '  -0.21%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.34'
$ws.Range('D44').ClearFormats()

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6560'
$ws.Range('D45').ClearFormats()

$ws.Range('E45').Value = '  -0.36%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.340'
$ws.Range('D46').ClearFormats()

$ws.Range('E46').Value = '  -0.84%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9991'
$ws.Range('D47').ClearFormats()

$ws.Range('E47').Value = '  +0.11%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.001'
$ws.Range('D48').ClearFormats()

$ws.Range('E48').Value = '  -1.07%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08045'
$ws.Range('D49').ClearFormats()

$ws.Range('E49').Value = '  +0.57%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.211'
$ws.Range('D50').ClearFormats()

$ws.Range('E50').Value = '  -0.19%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '125.54'
$ws.Range('D51').ClearFormats()

$ws.Range('E51').Value = '  -3.63%  '
